$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.235341333333333
$ws.Range("N2").Value = 9.706023999999999
$ws.Range("O2").Value = 0.2153734454473681
$ws.Range("P2").Value = 0.2153734454473681
$ws.Range("Q2").Value = 2.637541922937777
$ws.Range("R2").Value = 23.73787730644
$ws.Range("S2").Value = 0.2153734454473681
$ws.Range("T2").Value = 0.2153734454473681

# Row 3 updates
$ws.Range("O3").Value = 0.4841904166376352
$ws.Range("P3").Value = 0.4841904166376352
$ws.Range("S3").Value = 0.4841904166376352
$ws.Range("T3").Value = 0.4841904166376352

# Row 4 updates
$ws.Range("O4").Value = 0.3004361379149967
$ws.Range("P4").Value = 0.3004361379149967
$ws.Range("S4").Value = 0.3004361379149967
$ws.Range("T4").Value = 0.3004361379149967
